$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text in F1 ("First on Screen Appearance" -> "First Screen Appearance")
$ws.Range("F1").Value = "First Screen Appearance"

# Update the active selection to F1
$ws.Range("F1").Select()
